$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the register/search values (GaluaPulemet.. / Magnifikate..) in B2:C6
$ws.Range("B2").Value = "GaluaPulemet644"
$ws.Range("B3").Value = "GaluaPulemet724"
$ws.Range("B4").Value = "GaluaPulemet824"
$ws.Range("B5").Value = "GaluaPulemet924"
$ws.Range("B6").Value = "Magnifikate7724"

$ws.Range("C2").Value = "GaluaPulemet644@gmail.com"
$ws.Range("C3").Value = "GaluaPulemet724@gmail.com"
$ws.Range("C4").Value = "GaluaPulemet824@gmail.com"
$ws.Range("C5").Value = "GaluaPulemet924@gmail.com"
$ws.Range("C6").Value = "Magnifikate772@gmail.com"

# Update the active selection to K4
$ws.Range("K4").Select()
